$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol-list snapshot (Price / Volume(1h)) per the commit diff.
# Each value is written with a leading apostrophe so Excel stores it as
# literal text (matching the source workbook's inline-string cells) instead
# of coercing numeric-looking text (e.g. "306.26") or percent-looking text
# (e.g. "-0.06%") into a number. ClearFormats() then strips the "quote
# prefix" text-format style Excel stamps on such cells, so each cell keeps
# its original (default/no) style - exactly like the rest of the sheet's
# untouched text cells. ClearFormats() is called per-cell because multi-area
# (comma-separated) ranges only clear the first area.

$ws.Range("D2").Value = "'306.26"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-0.06%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'36.40"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-0.89%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.060"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'0.55%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07894"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'0.54%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'2.129"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-2.22%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'7.973"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.90%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.9272"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'0.35%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.09693"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-2.32%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1864"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'-0.86%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.09043"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.03717"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'2.79%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.09913"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.21%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001440"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-3.12%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.005620"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-1.03%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.468"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.19%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'4.157"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'2.33%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.660"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'13.95%"
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'-0.81%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.1313"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-2.57%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'5.112"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'3.57%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.2255"
$ws.Range("D22").ClearFormats()
$ws.Range("E23").Value = "'-0.93%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001239"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'0.47%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.004791"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'-7.16%"
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'74.07%"
$ws.Range("E27").ClearFormats()
$ws.Range("D39").Value = "'0.02007"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'10.84%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.04919"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'3.62%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007784"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-1.85%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.1394"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-0.85%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.007829"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'2.92%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.002145"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-1.81%"
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = "'11.42%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006294"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-0.71%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'0.00%"
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "'-0.03%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'51.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'42.98%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.001906"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'-29.21%"
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.00002104"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'0.00%"
$ws.Range("E51").ClearFormats()
